{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 0: title + paper name (line break between the two runs)\nparagraphs.items[0].insertText(\"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 05.03.25\\u000bMixtures of in-context learners\", Word.InsertLocation.replace);\n\n// Paragraph 1: intro text\nparagraphs.items[1].insertText(\"\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05de\u05d5\u05d3\u05e8\u05e0\u05d9\u05d9\u05dd \u05e0\u05d9\u05d7\u05e0\u05d9\u05dd \u05d1\u05d9\u05db\u05d5\u05dc\u05ea \u05dc\u05d1\u05e6\u05e2 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05d4\u05dd \u05dc\u05d0 \u05d0\u05d5\u05de\u05e0\u05d5 \u05e2\u05dc\u05d9\u05d4\u05dd \u05d1\u05d0\u05d5\u05e4\u05df \u05de\u05e4\u05d5\u05e8\u05e9 \u05d1\u05d4\u05ea\u05d1\u05e1\u05e1 \u05e2\u05dc \u05db\u05de\u05d4 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05d4\u05de\u05d3\u05d2\u05d9\u05de\u05d5\u05ea \u05d0\u05ea \u05d4\u05de\u05e9\u05d9\u05de\u05d4 \u05dc\u05dc\u05d0 \u05e6\u05d5\u05e8\u05da \u05d1\u05d0\u05d9\u05de\u05d5\u05df (\u05e4\u05d9\u05d9\u05df \u05d8\u05d9\u05d5\u05df).\u05d9\u05db\u05d5\u05dc\u05ea \u05d6\u05d5 \u05e7\u05d9\u05d1\u05dc\u05d4 \u05e9\u05dd \u05dc\u05de\u05d9\u05d3\u05d4 in-context (\u05d1\u05e7\u05e6\u05e8\u05d4 ICL) . \u05d0\u05e0\u05d9 \u05d2\u05dd \u05e8\u05d0\u05d9\u05ea\u05d9 \u05e9\u05e7\u05d5\u05e8\u05d0\u05d9\u05dd \u05dc\u05d6\u05d4 \u05dc\u05e4\u05e2\u05de\u05d9\u05dd \u05dc\u05de\u05d9\u05d3\u05ea few-shot \u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d6\u05d4 \u05e4\u05d7\u05d5\u05ea \u05de\u05ea\u05d0\u05d9\u05dd \u05db\u05d9 few-shot learning \u05de\u05d5\u05d2\u05d3\u05e8 \u05d1\u05d3\u05f4\u05db \u05d1\u05ea\u05d5\u05e8 \u05e4\u05d9\u05d9\u05df \u05d8\u05d9\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05e2\u05dc \u05db\u05de\u05d4 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea. \", Word.InsertLocation.replace);\n\n// Paragraph 2\nparagraphs.items[2].insertText(\"\u05d0\u05d6 \u05d0\u05d9\u05da \u05db\u05dc \u05d4\u05e2\u05e1\u05e7 \u05e2\u05d5\u05e1\u05e7? \u05de\u05e1\u05e4\u05e7\u05d9\u05dd \u05dc\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05db\u05de\u05d4 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e9\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2 \u05d4\u05de\u05e9\u05d9\u05de\u05d4 \u05d1\u05ea\u05d5\u05e8 \u05e4\u05e8\u05d5\u05de\u05e4\u05d8, \u05d1\u05d3\u05f4\u05db \u05db\u05de\u05d4 \u05d6\u05d5\u05d2\u05d5\u05ea \u05db\u05d0\u05e9\u05e8 x_i \u05d4\u05d9\u05e0\u05d4 \u05e9\u05d0\u05dc\u05d4 \u05d0\u05d5 \u05e9\u05d0\u05d9\u05dc\u05ea\u05d4\u05d4  \u05d5-y_i \u05d4\u05d9\u05e0\u05d4 \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05e6\u05e4\u05d5\u05d9\u05d4 \u05dc-x_i. \u05dc\u05d0\u05d7\u05e8 \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05d0\u05dc\u05d5 \u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 x \u05e9\u05d4\u05de\u05d5\u05d3\u05dc \u05e6\u05e8\u05d9\u05da \u05dc\u05e1\u05e4\u05e7 \u05ea\u05e9\u05d5\u05d1\u05d4 \u05e2\u05dc\u05d9\u05d4 \u05d1\u05d4\u05ea\u05d0\u05dd \u05dc\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05dc\u05e8\u05d0\u05d4 \u05dc\u05e4\u05e0\u05d9 \u05db\u05df. \", Word.InsertLocation.replace);\n\n// Paragraph 3\nparagraphs.items[3].insertText(\"\u05e1\u05d1\u05d9\u05e8 \u05dc\u05d4\u05e0\u05d9\u05d7 \u05dc\u05db\u05dc \u05e9\u05d0\u05dc\u05d4 x \u05d9\u05e9 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea x_i \u05d1\u05ea\u05d5\u05da \u05d4\u05e4\u05e8\u05d5\u05de\u05e4\u05d8 \u05e9\u05d3\u05d5\u05de\u05d5\u05ea \u05dc\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d5\u05d9\u05e9 \u05db\u05d0\u05dc\u05d5 \u05e9\u05e4\u05d7\u05d5\u05ea. \u05d0\u05d9\u05da \u05e0\u05d2\u05e8\u05d5\u05dd \u05dc\u05de\u05d5\u05d3\u05dc \u05dc\u05d4\u05ea\u05d7\u05e9\u05d1 \u05d9\u05d5\u05ea\u05e8 \u05d1\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d5\u05dc\u05d4\u05ea\u05d7\u05e9\u05d1 \u05e4\u05d7\u05d5\u05ea \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e4\u05d7\u05d5\u05ea \u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9\u05d5\u05ea \u05dc\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 x. \u05d6\u05d5 \u05d4\u05e9\u05d0\u05dc\u05d4 \u05e9\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05d5\u05d0\u05dc\u05d9\u05dd \u05d5\u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05d9\u05d8\u05d4 \u05dc\u05de\u05e9\u05e7\u05d5\u05dc \u05ea\u05e8\u05d5\u05de\u05d5\u05ea \u05e9\u05dc \u05db\u05dc \u05d3\u05d5\u05d2\u05de\u05d0 \u05dc\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4\u05d4 \u05e0\u05ea\u05d5\u05e0\u05d4 x. \", Word.InsertLocation.replace);\n\n// Paragraph 4\nparagraphs.items[4].insertText(\"\u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05dc \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05de\u05ea\u05d5\u05d9\u05d2\u05d5\u05ea (\u05e2\u05dd \u05ea\u05e9\u05d5\u05d1\u05d5\u05ea) \u05d5\u05de\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05d4\u05e4\u05d5\u05dc\u05d8 \u05de\u05e9\u05e7\u05dc w_i \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05d3\u05d5\u05d2\u05de\u05d0 \u05d1\u05e4\u05e8\u05d5\u05de\u05e4\u05d8 \u05dc\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 x. \u05de\u05e9\u05e7\u05d5\u05dc\u05d5\u05ea w_i \u05de\u05e9\u05de\u05e9\u05d5\u05ea \u05dc\u05d7\u05d9\u05e9\u05d5\u05d1 \u05e9\u05dc \u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e9\u05dc \u05db\u05dc \u05d8\u05d5\u05e7\u05df \u05d1\u05ea\u05e9\u05d5\u05d1\u05d4 y \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05db\u05dc \u05d6\u05d5\u05d2\u05d5\u05ea  \u05d5\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4\u05d4 x.  \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d4\u05d6\u05d5 \u05de\u05d9\u05d5\u05e6\u05d2\u05ea \u05d1\u05ea\u05d5\u05e8 \u05e1\u05db\u05d5\u05dd \u05de\u05de\u05d5\u05e9\u05e7\u05dc \u05e2\u05dd w_I \u05e9\u05dc log-probs \u05e9\u05dc \u05d8\u05d5\u05e7\u05df y \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05db\u05dc \u05d6\u05d5\u05d2 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea . \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05e9\u05ea\u05d9 \u05d3\u05e8\u05db\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05d0\u05ea \u05d4\u05de\u05e9\u05e7\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05dc\u05d5 (\u05e2\u05dc \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05dc \u05e9\u05d0\u05dc\u05d5\u05ea \u05d5\u05ea\u05e9\u05d5\u05d1\u05d5\u05ea). \u05d4\u05d3\u05e8\u05da \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4 \u05dc\u05d0\u05de\u05df \u05d0\u05d5\u05ea\u05d4 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d9\u05e9\u05d9\u05e8\u05d4 (\u05e4\u05e9\u05d5\u05d8 \u05dc\u05d0\u05e4\u05d8\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05ea \u05dc\u05d5\u05e1 \u05dc\u05e4\u05d9\u05d4\u05df \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e9\u05dc x_i \u05d5-y_i) \u05d5\u05d4\u05d3\u05e8\u05da \u05d4\u05e9\u05e0\u05d9\u05d9\u05d4 \u05d4\u05d9\u05d0 \u05dc\u05d0\u05de\u05df \u05e8\u05e9\u05ea \u05d4\u05de\u05d7\u05e9\u05d1\u05ea \u05d0\u05ea \u05d4\u05de\u05e9\u05e7\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05dc\u05d5 \u05d5\u05dc\u05d0\u05e4\u05d8\u05dd \u05d0\u05ea \u05d4\u05de\u05e9\u05e7\u05d5\u05dc\u05d5\u05ea \u05e9\u05dc\u05d4.\", Word.InsertLocation.replace);\n\n// Paragraph 5\nparagraphs.items[5].insertText(\"\u05d1\u05e1\u05d5\u05e3 \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05d9\u05d8\u05d4 \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc top-k \u05e9\u05dc \u05d4\u05de\u05e9\u05e7\u05dc\u05d9\u05dd \u05db\u05d3\u05d9 \u05dc\u05d0 \u05dc\u05d7\u05e9\u05d1 \u05d0\u05ea \u05db\u05dc \u05d4-log-probs \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e9\u05d6\u05d4 \u05d9\u05db\u05d5\u05dc \u05dc\u05d4\u05d9\u05d5\u05ea \u05e7\u05e6\u05ea \u05db\u05d1\u05d3 \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea \u05d5\u05d2\u05dd \u05dc\u05d5\u05e7\u05d7 \u05d6\u05de\u05df. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e2\u05dc Implicit MLE \u05e9\u05d4\u05d9\u05d0 \u05de\u05d0\u05de\u05e0\u05ea \u05de\u05d5\u05d3\u05dc \u05dc\u05d0\u05e4\u05d8\u05dd \u05de\u05d5\u05d3\u05dc \u05dc\u05d8\u05e0\u05d8\u05d9 \u05db\u05d0\u05e9\u05e8 \u05de\u05e9\u05ea\u05e0\u05d5 \u05d7\u05d1\u05d5\u05d9 (\u05dc\u05d8\u05e0\u05d8\u05d9) \u05e0\u05d3\u05d2\u05dd \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d3\u05d9\u05e1\u05e7\u05e8\u05d8\u05d9\u05ea. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d3\u05d9 \u05dc\u05d0 \u05d8\u05e8\u05d9\u05d5\u05d5\u05d9\u05d0\u05dc\u05d9\u05ea \u05dc\u05d4\u05d1\u05e0\u05d4 - \u05de\u05d9 \u05e9\u05e8\u05d5\u05e6\u05d4 \u05dc\u05d4\u05ea\u05e2\u05de\u05e7 \u05d1\u05d4 (\u05de\u05d5\u05de\u05dc\u05e5) \u05de\u05d5\u05d6\u05de\u05df \u05dc\u05d4\u05d1\u05d9\u05d8 \u05d1\u05e8\u05e4\u05e8\u05e0\u05e1\u05d9\u05dd.\", Word.InsertLocation.replace);\n\nawait context.sync();\n\n// Remove the ten paragraphs (indices 6-15) that described the old paper's\n// bullet-point observations -- delete from the end so earlier indices stay valid.\nfor (let i = 15; i >= 6; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\n// Paragraph 16 (now at index 6 after the deletions): the arxiv link\nconst remaining = body.paragraphs;\nremaining.load(\"items\");\nawait context.sync();\nconst linkParagraph = remaining.items[remaining.items.length - 1];\nlinkParagraph.insertText(\"https://arxiv.org/abs/2411.02830\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Paragraph 1: title + paper name (line break between the two runs)\n$d.Paragraphs.Item(1).Range.Text = \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 05.03.25\u000bMixtures of in-context learners\"\n\n# Paragraph 2: intro text\n$d.Paragraphs.Item(2).Range.Text = \"\u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05de\u05d5\u05d3\u05e8\u05e0\u05d9\u05d9\u05dd \u05e0\u05d9\u05d7\u05e0\u05d9\u05dd \u05d1\u05d9\u05db\u05d5\u05dc\u05ea \u05dc\u05d1\u05e6\u05e2 \u05de\u05e9\u05d9\u05de\u05d5\u05ea \u05e9\u05d4\u05dd \u05dc\u05d0 \u05d0\u05d5\u05de\u05e0\u05d5 \u05e2\u05dc\u05d9\u05d4\u05dd \u05d1\u05d0\u05d5\u05e4\u05df \u05de\u05e4\u05d5\u05e8\u05e9 \u05d1\u05d4\u05ea\u05d1\u05e1\u05e1 \u05e2\u05dc \u05db\u05de\u05d4 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05d4\u05de\u05d3\u05d2\u05d9\u05de\u05d5\u05ea \u05d0\u05ea \u05d4\u05de\u05e9\u05d9\u05de\u05d4 \u05dc\u05dc\u05d0 \u05e6\u05d5\u05e8\u05da \u05d1\u05d0\u05d9\u05de\u05d5\u05df (\u05e4\u05d9\u05d9\u05df \u05d8\u05d9\u05d5\u05df).\u05d9\u05db\u05d5\u05dc\u05ea \u05d6\u05d5 \u05e7\u05d9\u05d1\u05dc\u05d4 \u05e9\u05dd \u05dc\u05de\u05d9\u05d3\u05d4 in-context (\u05d1\u05e7\u05e6\u05e8\u05d4 ICL) . \u05d0\u05e0\u05d9 \u05d2\u05dd \u05e8\u05d0\u05d9\u05ea\u05d9 \u05e9\u05e7\u05d5\u05e8\u05d0\u05d9\u05dd \u05dc\u05d6\u05d4 \u05dc\u05e4\u05e2\u05de\u05d9\u05dd \u05dc\u05de\u05d9\u05d3\u05ea few-shot \u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d6\u05d4 \u05e4\u05d7\u05d5\u05ea \u05de\u05ea\u05d0\u05d9\u05dd \u05db\u05d9 few-shot learning \u05de\u05d5\u05d2\u05d3\u05e8 \u05d1\u05d3\u05f4\u05db \u05d1\u05ea\u05d5\u05e8 \u05e4\u05d9\u05d9\u05df \u05d8\u05d9\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05e2\u05dc \u05db\u05de\u05d4 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea. \"\n\n# Paragraph 3\n$d.Paragraphs.Item(3).Range.Text = \"\u05d0\u05d6 \u05d0\u05d9\u05da \u05db\u05dc \u05d4\u05e2\u05e1\u05e7 \u05e2\u05d5\u05e1\u05e7? \u05de\u05e1\u05e4\u05e7\u05d9\u05dd \u05dc\u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05db\u05de\u05d4 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e9\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2 \u05d4\u05de\u05e9\u05d9\u05de\u05d4 \u05d1\u05ea\u05d5\u05e8 \u05e4\u05e8\u05d5\u05de\u05e4\u05d8, \u05d1\u05d3\u05f4\u05db \u05db\u05de\u05d4 \u05d6\u05d5\u05d2\u05d5\u05ea \u05db\u05d0\u05e9\u05e8 x_i \u05d4\u05d9\u05e0\u05d4 \u05e9\u05d0\u05dc\u05d4 \u05d0\u05d5 \u05e9\u05d0\u05d9\u05dc\u05ea\u05d4\u05d4  \u05d5-y_i \u05d4\u05d9\u05e0\u05d4 \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05d4\u05e6\u05e4\u05d5\u05d9\u05d4 \u05dc-x_i. \u05dc\u05d0\u05d7\u05e8 \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05d0\u05dc\u05d5 \u05de\u05d6\u05d9\u05e0\u05d9\u05dd \u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 x \u05e9\u05d4\u05de\u05d5\u05d3\u05dc \u05e6\u05e8\u05d9\u05da \u05dc\u05e1\u05e4\u05e7 \u05ea\u05e9\u05d5\u05d1\u05d4 \u05e2\u05dc\u05d9\u05d4 \u05d1\u05d4\u05ea\u05d0\u05dd \u05dc\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05dc\u05e8\u05d0\u05d4 \u05dc\u05e4\u05e0\u05d9 \u05db\u05df. \"\n\n# Paragraph 4\n$d.Paragraphs.Item(4).Range.Text = \"\u05e1\u05d1\u05d9\u05e8 \u05dc\u05d4\u05e0\u05d9\u05d7 \u05dc\u05db\u05dc \u05e9\u05d0\u05dc\u05d4 x \u05d9\u05e9 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea x_i \u05d1\u05ea\u05d5\u05da \u05d4\u05e4\u05e8\u05d5\u05de\u05e4\u05d8 \u05e9\u05d3\u05d5\u05de\u05d5\u05ea \u05dc\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d5\u05d9\u05e9 \u05db\u05d0\u05dc\u05d5 \u05e9\u05e4\u05d7\u05d5\u05ea. \u05d0\u05d9\u05da \u05e0\u05d2\u05e8\u05d5\u05dd \u05dc\u05de\u05d5\u05d3\u05dc \u05dc\u05d4\u05ea\u05d7\u05e9\u05d1 \u05d9\u05d5\u05ea\u05e8 \u05d1\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d5\u05dc\u05d4\u05ea\u05d7\u05e9\u05d1 \u05e4\u05d7\u05d5\u05ea \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e4\u05d7\u05d5\u05ea \u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9\u05d5\u05ea \u05dc\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 x. \u05d6\u05d5 \u05d4\u05e9\u05d0\u05dc\u05d4 \u05e9\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d4\u05de\u05d0\u05de\u05e8 \u05e9\u05d5\u05d0\u05dc\u05d9\u05dd \u05d5\u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05d9\u05d8\u05d4 \u05dc\u05de\u05e9\u05e7\u05d5\u05dc \u05ea\u05e8\u05d5\u05de\u05d5\u05ea \u05e9\u05dc \u05db\u05dc \u05d3\u05d5\u05d2\u05de\u05d0 \u05dc\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4\u05d4 \u05e0\u05ea\u05d5\u05e0\u05d4 x. \"\n\n# Paragraph 5\n$d.Paragraphs.Item(5).Range.Text = \"\u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05dc \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05de\u05ea\u05d5\u05d9\u05d2\u05d5\u05ea (\u05e2\u05dd \u05ea\u05e9\u05d5\u05d1\u05d5\u05ea) \u05d5\u05de\u05d0\u05de\u05df \u05de\u05d5\u05d3\u05dc \u05d4\u05e4\u05d5\u05dc\u05d8 \u05de\u05e9\u05e7\u05dc w_i \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05d3\u05d5\u05d2\u05de\u05d0 \u05d1\u05e4\u05e8\u05d5\u05de\u05e4\u05d8 \u05dc\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4 x. \u05de\u05e9\u05e7\u05d5\u05dc\u05d5\u05ea w_i \u05de\u05e9\u05de\u05e9\u05d5\u05ea \u05dc\u05d7\u05d9\u05e9\u05d5\u05d1 \u05e9\u05dc \u05d4\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05e9\u05dc \u05db\u05dc \u05d8\u05d5\u05e7\u05df \u05d1\u05ea\u05e9\u05d5\u05d1\u05d4 y \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05db\u05dc \u05d6\u05d5\u05d2\u05d5\u05ea  \u05d5\u05e9\u05d0\u05d9\u05dc\u05ea\u05d4\u05d4 x.  \u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d4\u05d6\u05d5 \u05de\u05d9\u05d5\u05e6\u05d2\u05ea \u05d1\u05ea\u05d5\u05e8 \u05e1\u05db\u05d5\u05dd \u05de\u05de\u05d5\u05e9\u05e7\u05dc \u05e2\u05dd w_I \u05e9\u05dc log-probs \u05e9\u05dc \u05d8\u05d5\u05e7\u05df y \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05db\u05dc \u05d6\u05d5\u05d2 \u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea . \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05e2 \u05e9\u05ea\u05d9 \u05d3\u05e8\u05db\u05d9\u05dd \u05dc\u05d0\u05de\u05df \u05d0\u05ea \u05d4\u05de\u05e9\u05e7\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05dc\u05d5 (\u05e2\u05dc \u05d3\u05d0\u05d8\u05d4\u05e1\u05d8 \u05e9\u05dc \u05e9\u05d0\u05dc\u05d5\u05ea \u05d5\u05ea\u05e9\u05d5\u05d1\u05d5\u05ea). \u05d4\u05d3\u05e8\u05da \u05d4\u05e8\u05d0\u05e9\u05d5\u05e0\u05d4 \u05dc\u05d0\u05de\u05df \u05d0\u05d5\u05ea\u05d4 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d9\u05e9\u05d9\u05e8\u05d4 (\u05e4\u05e9\u05d5\u05d8 \u05dc\u05d0\u05e4\u05d8\u05dd \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05ea \u05dc\u05d5\u05e1 \u05dc\u05e4\u05d9\u05d4\u05df \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e9\u05dc x_i \u05d5-y_i) \u05d5\u05d4\u05d3\u05e8\u05da \u05d4\u05e9\u05e0\u05d9\u05d9\u05d4 \u05d4\u05d9\u05d0 \u05dc\u05d0\u05de\u05df \u05e8\u05e9\u05ea \u05d4\u05de\u05d7\u05e9\u05d1\u05ea \u05d0\u05ea \u05d4\u05de\u05e9\u05e7\u05d5\u05dc\u05d5\u05ea \u05d4\u05d0\u05dc\u05d5 \u05d5\u05dc\u05d0\u05e4\u05d8\u05dd \u05d0\u05ea \u05d4\u05de\u05e9\u05e7\u05d5\u05dc\u05d5\u05ea \u05e9\u05dc\u05d4.\"\n\n# Paragraph 6\n$d.Paragraphs.Item(6).Range.Text = \"\u05d1\u05e1\u05d5\u05e3 \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05e9\u05d9\u05d8\u05d4 \u05dc\u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc top-k \u05e9\u05dc \u05d4\u05de\u05e9\u05e7\u05dc\u05d9\u05dd \u05db\u05d3\u05d9 \u05dc\u05d0 \u05dc\u05d7\u05e9\u05d1 \u05d0\u05ea \u05db\u05dc \u05d4-log-probs \u05e2\u05d1\u05d5\u05e8 \u05db\u05dc \u05d4\u05d3\u05d5\u05d2\u05de\u05d0\u05d5\u05ea \u05e9\u05d6\u05d4 \u05d9\u05db\u05d5\u05dc \u05dc\u05d4\u05d9\u05d5\u05ea \u05e7\u05e6\u05ea \u05db\u05d1\u05d3 \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea \u05d5\u05d2\u05dd \u05dc\u05d5\u05e7\u05d7 \u05d6\u05de\u05df. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05de\u05d1\u05d5\u05e1\u05e1\u05ea \u05e2\u05dc Implicit MLE \u05e9\u05d4\u05d9\u05d0 \u05de\u05d0\u05de\u05e0\u05ea \u05de\u05d5\u05d3\u05dc \u05dc\u05d0\u05e4\u05d8\u05dd \u05de\u05d5\u05d3\u05dc \u05dc\u05d8\u05e0\u05d8\u05d9 \u05db\u05d0\u05e9\u05e8 \u05de\u05e9\u05ea\u05e0\u05d5 \u05d7\u05d1\u05d5\u05d9 (\u05dc\u05d8\u05e0\u05d8\u05d9) \u05e0\u05d3\u05d2\u05dd \u05de\u05d4\u05ea\u05e4\u05dc\u05d2\u05d5\u05ea \u05d3\u05d9\u05e1\u05e7\u05e8\u05d8\u05d9\u05ea. \u05d4\u05e9\u05d9\u05d8\u05d4 \u05d3\u05d9 \u05dc\u05d0 \u05d8\u05e8\u05d9\u05d5\u05d5\u05d9\u05d0\u05dc\u05d9\u05ea \u05dc\u05d4\u05d1\u05e0\u05d4 - \u05de\u05d9 \u05e9\u05e8\u05d5\u05e6\u05d4 \u05dc\u05d4\u05ea\u05e2\u05de\u05e7 \u05d1\u05d4 (\u05de\u05d5\u05de\u05dc\u05e5) \u05de\u05d5\u05d6\u05de\u05df \u05dc\u05d4\u05d1\u05d9\u05d8 \u05d1\u05e8\u05e4\u05e8\u05e0\u05e1\u05d9\u05dd.\"\n\n# Remove the ten paragraphs (7-16, 1-indexed) that described the old paper's\n# bullet-point observations -- delete from the end so earlier indices stay valid.\nfor ($i = 16; $i -ge 7; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n\n# Paragraph 17 (now the last paragraph after the deletions): the arxiv link\n$last = $d.Paragraphs.Item($d.Paragraphs.Count)\n$last.Range.Text = \"https://arxiv.org/abs/2411.02830\"\n\n"}
